$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: split a run in three by locating "prefix<mid>suffix" text
# and nudging a direct-formatting property on the <mid> sub-range so
# the engine has to materialise it as its own <w:r>. The toggle is
# applied and then reverted so the rendered formatting is unchanged.
# `beforeText` is the (unique) text immediately preceding `midText`
# within the match, used to disambiguate repeated substrings.
# ------------------------------------------------------------------
function Split-MiddleRun($searchText, $beforeText, $midText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $searchText)
        return $null
    }
    $whole = $find.Parent
    $wholeStart = $whole.Start
    $start = $wholeStart + $beforeText.Length
    $end = $start + $midText.Length
    $mid = $d.Range($start, $end)
    if ($mid.Text -ne $midText) {
        Write-Output ("MISMATCH: expected [" + $midText + "] got [" + $mid.Text + "]")
    }
    $mid.Bold = 1
    $mid.Bold = 0
    return $mid
}

# ------------------------------------------------------------------
# Helper: replace `oldSearch` with `newReplace` (a plain text swap,
# e.g. fixing a missing accent) and, once the new text is in place,
# force the `newMid` word inside it to live in its own <w:r> (mirrors
# the source diff, which always isolates the touched word into a
# separate run even when no proofErr wrapper is involved).
# ------------------------------------------------------------------
function Replace-AndSplitMiddle($oldSearch, $newReplace, $beforeText, $newMid) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $ok = $find.Execute($oldSearch, $true, $false, $false, $false, $false, $true, 1, $false, $newReplace, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $oldSearch)
        return $null
    }
    $find2 = $d.Content.Find
    $find2.ClearFormatting()
    $ok2 = $find2.Execute($newReplace, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok2) {
        Write-Output ("NOT FOUND (pass2): " + $newReplace)
        return $null
    }
    $whole = $find2.Parent
    $wholeStart = $whole.Start
    $start = $wholeStart + $beforeText.Length
    $end = $start + $newMid.Length
    $mid = $d.Range($start, $end)
    if ($mid.Text -ne $newMid) {
        Write-Output ("MISMATCH: expected [" + $newMid + "] got [" + $mid.Text + "]")
    }
    $mid.Bold = 1
    $mid.Bold = 0
    return $mid
}

# 1) "Controle de Estoque" description: cosmetic spell-check split
#    around "a" -- text itself is unchanged.
Split-MiddleRun "Deixando claro todos os itens que tem a disposição" "Deixando claro todos os itens que tem " "a" | Out-Null

# 2) "Cálculo de comissões" description: cosmetic spell-check split
#    around "pré" -- text itself is unchanged.
Split-MiddleRun "Fazendo os cálculos pré-configurados" "Fazendo os cálculos " "pré" | Out-Null

# 3) Row 14 title: Funcionalidades -> Divulgação em vídeo
$d.Content.Find.Execute("Funcionalidades", $true, $false, $false, $false, $false, $true, 1, $false, "Divulgação em vídeo", 2) | Out-Null

# 4) Row 15 title: Comunicação -> Comunicação com cliente (appended as
#    its own trailing run, the original "Comunicação" run untouched)
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Comunicação", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r = $find.Parent
    $insStart = $r.End
    $r.Collapse(0)
    $suffix = " com cliente"
    $r.InsertAfter($suffix)
    $newRange = $d.Range($insStart, $insStart + $suffix.Length)
    $newRange.Bold = 1
    $newRange.Bold = 0
}

# 5) Row 15 description: "numero" -> "número" (real spelling fix),
#    isolated into its own run.
Replace-AndSplitMiddle "Deixaremos disponível também um numero telefônico" "Deixaremos disponível também um número telefônico" "Deixaremos disponível também um " "número" | Out-Null

# 6) Row 18 (Newsletter) description: cosmetic spell-check split
#    around "emails" -- text itself is unchanged.
Split-MiddleRun "Esses emails serão disparados" "Esses " "emails" | Out-Null

# 7) Row 19 (Perfil de clientes) description: "publico" -> "público",
#    isolated into its own run.
Replace-AndSplitMiddle "entenda melhor seu publico, mostrando" "entenda melhor seu público, mostrando" "entenda melhor seu " "público" | Out-Null

# 8) Row 25 (Histórico de atendimento) description: "criticas" -> "críticas",
#    isolated into its own run.
Replace-AndSplitMiddle "sugestões e criticas feitas pelos clientes" "sugestões e críticas feitas pelos clientes" "sugestões e " "críticas" | Out-Null

# 9) Delete the last table row (26 / Cadastro cliente)
$t = $d.Tables.Item(1)
$t.Rows.Item($t.Rows.Count).Delete()
